$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: SetHighscore response text changed + row grew taller ---
$ws.Range("F6").Value = "OK 200 (No error, but old value was higher),`nCREATED 201 (New highscore stored)`nor`nerrorCode"
$ws.Rows.Item(6).RowHeight = 120

# --- Column C got wider to fit the new, longer route string ---
$ws.Columns.Item(3).ColumnWidth = 48.75

# --- New row 8: GetRanking endpoint ---
$ws.Range("A8").Value = "GetRanking"
$ws.Range("B8").Value = "Get"
$ws.Range("C8").Value = "/ranking/{levelIndex}[?limit={maxNumberOfResults}]"
$ws.Range("D8").Value = "no"
$ws.Range("F8").Value = "OK 200,`nor`nerrorCode"
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("F8").VerticalAlignment = -4108
$ws.Range("F8").WrapText = $true
$ws.Range("G8").Value = "{`n<number (index)> : {`nusername: <string>,`nvalue: <number>`n}`n}"
$ws.Range("G8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 90

# --- View state: scrolled down, D12 selected ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()
